# Generate Report for Handback
# Updates the "Overview" sheet status text, and records an error detail
# message (with new column width) on the per-language handback report
# sheets "zh-cn" and "de-de" for the rows where the handback file name
# did not match the expected handoff file name.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status changed from "Ready for handoff" to
#     "Handback transform failed" for the b8d1e892 file row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# Column width (in ColumnWidth units) that renders as an OOXML column
# width of exactly 40 once Excel applies its fixed padding offset.
$errorDetailColumnWidth = 39.166666666666664

# --- zh-cn sheet: the "Status" cell for the same file row also reads
#     "Ready for handoff" and must be updated, and the "Error Detail"
#     column (P) for row 3 gets the transform failure message, with the
#     column widened so the message is readable.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth
$wsZhCn.Range("P3").Value = "Handback file name: 1eiw2fjm.2vc is different with handoff file name: b8d1e892-c2b9-481e-92b4-b50a62ac49fd.5167ad2236a28e0b64653a7e23cf6187c122eb1d.zh-cn."

# --- de-de sheet: same treatment as zh-cn.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth
$wsDeDe.Range("P3").Value = "Handback file name: 1eiw2fjm.2vc is different with handoff file name: b8d1e892-c2b9-481e-92b4-b50a62ac49fd.5167ad2236a28e0b64653a7e23cf6187c122eb1d.de-de."
